# Apply updated values to cfs_3_0.9 sheet as per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25
$ws.Range("I25").Value = 248889

# Row 28
$ws.Range("R28").Value = 79
$ws.Range("W28").Value = 3

# Row 33
$ws.Range("I33").Value = 248889

# Row 34
$ws.Range("I34").Value = 253918

# Row 35
$ws.Range("I35").Value = 248938

# Row 36
$ws.Range("U36").Value = 79

# Row 39
$ws.Range("C39").Value = 85
$ws.Range("D39").Value = 85
$ws.Range("F39").Value = 6
$ws.Range("G39").Value = 30.57
$ws.Range("H39").Value = 24.74
$ws.Range("P39").Value = 2
$ws.Range("U39").Value = 79

# Row 40
$ws.Range("U40").Value = 73

# Row 48
$ws.Range("U48").Value = 73

# Row 49
$ws.Range("U49").Value = 73

# Row 51
$ws.Range("U51").Value = 73

# Row 53
$ws.Range("U53").Value = 73

# Row 54
$ws.Range("U54").Value = 73

# Row 55
$ws.Range("U55").Value = 73

# Row 56
$ws.Range("U56").Value = 73

# Row 58
$ws.Range("C58").Value = 85
$ws.Range("D58").Value = 85
$ws.Range("F58").Value = 6
$ws.Range("G58").Value = 30.57
$ws.Range("H58").Value = 24.54
$ws.Range("U58").Value = 79

# Row 62
$ws.Range("U62").Value = 73

# Row 65
$ws.Range("U65").Value = 73

# Row 68
$ws.Range("U68").Value = 79

# Row 70
$ws.Range("U70").Value = 73
